$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" column header in H1, reusing the exact same formatting
# (bold, bordered, centered header style) as the existing header cells by
# copying the format from the adjacent "sum" header (G1) instead of
# re-building it property-by-property.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Populate the new Save column with 0 for each existing data row
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
